$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 79.666664
$ws.Range("I9").Value = 70.875
$ws.Range("K9").Value = 70.875
$ws.Range("M9").Value = 98.125
$ws.Range("H15").Value = 849.9677
$ws.Range("I15").Value = 849.9677
$ws.Range("K15").Value = 2549.9031
$ws.Range("M15").Value = -2380.9031
$ws.Range("H33").Value = 168.28572
$ws.Range("I33").Value = 146.33333
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 146.33333
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = 82.66667000000001
$ws.Range("N33").Value = -758
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4350
$ws.Range("H70").Value = 2631.3333
$ws.Range("J70").Value = 2631.3333
$ws.Range("L70").Value = 7893.999899999999
$ws.Range("N70").Value = -8433.999899999999
$ws.Range("H73").Value = 2631.3333
$ws.Range("J73").Value = 2631.3333
$ws.Range("L73").Value = 7893.999899999999
$ws.Range("N73").Value = -9765.999899999999
$ws.Range("H115").Value = 9185
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H138").Value = 50010000
$ws.Range("I138").Value = 50010000
$ws.Range("K138").Value = 150030000
$ws.Range("M138").Value = -150024860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5583.3335
$ws.Range("I26").Value = 5468.75
$ws.Range("J26").Value = 6500
$ws.Range("K26").Value = 5468.75
$ws.Range("L26").Value = 6500
$ws.Range("M26").Value = -5138.75
$ws.Range("N26").Value = -7160
$ws.Range("H98").Value = 61977.5
$ws.Range("J98").Value = 61977.5
$ws.Range("L98").Value = 61977.5
$ws.Range("N98").Value = -67967.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 175
$ws.Range("I8").Value = 100
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 100
$ws.Range("L8").Value = 250
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = -530
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("N17").Value = -844
$ws.Range("H86").Value = 1750
$ws.Range("I86").Value = 1750
$ws.Range("K86").Value = 1750
$ws.Range("M86").Value = -627
$ws.Range("H89").Value = 1750
$ws.Range("I89").Value = 1750
$ws.Range("K89").Value = 8750
$ws.Range("M89").Value = -3134
$ws.Range("H134").Value = 6057.625
$ws.Range("I134").Value = 3942.2
$ws.Range("K134").Value = 11826.6
$ws.Range("M134").Value = -9291.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H58").Value = 6437
$ws.Range("I58").Value = 3546.25
$ws.Range("K58").Value = 3546.25
$ws.Range("M58").Value = -3343.25
$ws.Range("H103").Value = 1445996.2
$ws.Range("I103").Value = 20329
$ws.Range("J103").Value = 10000000
$ws.Range("K103").Value = 20329
$ws.Range("L103").Value = 10000000
$ws.Range("M103").Value = -19157
$ws.Range("N103").Value = -10002344
$ws.Range("H132").Value = 4636.522
$ws.Range("J132").Value = 7068
$ws.Range("L132").Value = 21204
$ws.Range("N132").Value = -26264
$ws.Range("H134").Value = 5188.533
$ws.Range("I134").Value = 2644.5557
$ws.Range("K134").Value = 7933.6671
$ws.Range("M134").Value = -5398.6671
$ws.Range("H136").Value = 6437
$ws.Range("I136").Value = 3546.25
$ws.Range("K136").Value = 10638.75
$ws.Range("M136").Value = -8088.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 123500
$ws.Range("J15").Value = 123500
$ws.Range("L15").Value = 123500
$ws.Range("N15").Value = -124076
$ws.Range("H81").Value = 123500
$ws.Range("J81").Value = 123500
$ws.Range("L81").Value = 123500
$ws.Range("N81").Value = -125496
$ws.Range("H84").Value = 123500
$ws.Range("J84").Value = 123500
$ws.Range("L84").Value = 370500
$ws.Range("N84").Value = -380484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 12800
$ws.Range("J3").Value = 12800
$ws.Range("L3").Value = 12800
$ws.Range("N3").Value = -13024
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H14").Value = 12800
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 12800
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 12800
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -13144
$ws.Range("H15").Value = 12800
$ws.Range("J15").Value = 12800
$ws.Range("L15").Value = 12800
$ws.Range("N15").Value = -13140
$ws.Range("H22").Value = 2826.25
$ws.Range("I22").Value = 2152.5
$ws.Range("J22").Value = 3500
$ws.Range("K22").Value = 2152.5
$ws.Range("L22").Value = 3500
$ws.Range("M22").Value = -1857.5
$ws.Range("N22").Value = -4090
$ws.Range("H24").Value = 19900
$ws.Range("J24").Value = 19900
$ws.Range("L24").Value = 19900
$ws.Range("N24").Value = -20586
$ws.Range("H27").Value = 2826.25
$ws.Range("I27").Value = 2152.5
$ws.Range("J27").Value = 3500
$ws.Range("K27").Value = 2152.5
$ws.Range("L27").Value = 3500
$ws.Range("M27").Value = -2045.5
$ws.Range("N27").Value = -3714
$ws.Range("H55").Value = 437.625
$ws.Range("I55").Value = 416.83334
$ws.Range("K55").Value = 416.83334
$ws.Range("M55").Value = -243.83334
$ws.Range("H68").Value = 3200
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3200
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H125").Value = 19999
$ws.Range("J125").Value = 19999
$ws.Range("L125").Value = 19999
$ws.Range("N125").Value = -29839
$ws.Range("H132").Value = 15781.857
$ws.Range("I132").Value = 15412.167
$ws.Range("K132").Value = 46236.501
$ws.Range("M132").Value = -43706.501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 40750
$ws.Range("I2").Value = 31500
$ws.Range("J2").Value = 50000
$ws.Range("K2").Value = 31500
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = -31388
$ws.Range("N2").Value = -50224
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("K14").Value = 500
$ws.Range("M14").Value = -332
$ws.Range("H18").Value = 36600
$ws.Range("J18").Value = 19900
$ws.Range("L18").Value = 19900
$ws.Range("N18").Value = -20246
$ws.Range("H20").Value = 22600
$ws.Range("J20").Value = 22600
$ws.Range("L20").Value = 22600
$ws.Range("N20").Value = -23080
$ws.Range("H81").Value = 401
$ws.Range("I81").Value = 401
$ws.Range("K81").Value = 802
$ws.Range("M81").Value = 259
$ws.Range("H84").Value = 401
$ws.Range("I84").Value = 401
$ws.Range("K84").Value = 4010
$ws.Range("M84").Value = 1294
